{"js": "// Fix grammar in the MS4 Scrum report:\n//  1. \"Review edthe purpose of white box tests and assigned each member\n//      whitebox test(s) to write\"\n//     -> \"Reviewed the purpose of white box tests and assigned each member\n//         white box test(s) to write\"\n//  2. \"Review matrix requirements again and discuss what needs to be revised\"\n//     -> \"Reviewed matrix requirements again and discuss what needs to be\n//         revised\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst OLD_TEXT_1 =\n  \"Review edthe purpose of white box tests and assigned each member whitebox test(s) to write\";\nconst NEW_TEXT_1 =\n  \"Reviewed the purpose of white box tests and assigned each member white box test(s) to write\";\n\nconst OLD_TEXT_2 =\n  \"Review matrix requirements again and discuss what needs to be revised\";\nconst NEW_TEXT_2 =\n  \"Reviewed matrix requirements again and discuss what needs to be revised\";\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (text === OLD_TEXT_1) {\n    paragraph.getRange().insertText(NEW_TEXT_1, \"Replace\");\n  } else if (text.trim() === OLD_TEXT_2) {\n    paragraph.getRange().insertText(NEW_TEXT_2, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix grammar in the MS4 Scrum report:\n#  1. \"Review edthe purpose of white box tests and assigned each member\n#      whitebox test(s) to write\"\n#     -> \"Reviewed the purpose of white box tests and assigned each member\n#         white box test(s) to write\"\n#  2. \"Review matrix requirements again and discuss what needs to be revised\"\n#     -> \"Reviewed matrix requirements again and discuss what needs to be\n#         revised\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Replacement.ClearFormatting()\n$range1.Find.Execute(\n    \"Review edthe purpose of white box tests and assigned each member whitebox test(s) to write\",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Reviewed the purpose of white box tests and assigned each member white box test(s) to write\",\n    $wdReplaceAll\n)\n\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Replacement.ClearFormatting()\n$range2.Find.Execute(\n    \"Review matrix requirements again and discuss what needs to be revised\",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Reviewed matrix requirements again and discuss what needs to be revised\",\n    $wdReplaceAll\n)\n"}
